$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# per the "Updated cryptos list" GitHub Actions commit.
$ws.Range("D2").Value = '37.266.89'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '2.069.30'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'233.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'56.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").Value = "'0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").Value = "'0.0764"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").Value = '2.372.61'
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").Value = "'14.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.87%  '
$ws.Range("D14").Value = "'20.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").Value = "'0.777"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").Value = '2.067.05'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '37.262.25'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = "'6.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.37%  '
$ws.Range("D20").Value = "'69.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = '0.0₃0814'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = "'226.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = "'166.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -3.84%  '
$ws.Range("D31").Value = "'0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = "'4.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = "'0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.12%  '
$ws.Range("D34").Value = "'4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("E35").Value = '  -4.66%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = "'1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("E38").Value = '  -3.48%  '
$ws.Range("E39").Value = '  -3.68%  '
$ws.Range("D40").Value = "'2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").Value = "'4.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.06%  '
$ws.Range("D42").Value = '1.481.15'
$ws.Range("D43").Value = "'96.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").Value = "'0.0932"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = "'15.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.97%  '
$ws.Range("D49").Value = "'7.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").Value = '2.260.78'
$ws.Range("E51").Value = '  -0.43%  '
